$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns C through J
$ws.Range("C1").Value = "Onderwerp"
$ws.Range("D1").Value = "Afzender"
$ws.Range("E1").Value = "Categorie"
$ws.Range("F1").Value = "Tijdstip"
$ws.Range("G1").Value = "Beantwoord"
$ws.Range("H1").Value = "Handmatig opvolgen"
$ws.Range("I1").Value = "Automatisch afgehandeld"
$ws.Range("J1").Value = "Hergebruikt antwoord"

# Copy the header style from B1 to the new header cells
$ws.Range("B1").Copy()
$ws.Range("C1:J1").PasteSpecial(-4122)

# Append new row 20 with full data across all columns
$ws.Range("A20").Value = "Testmail #4: Wil je 100 stuks M5-bouten bestellen?"
$ws.Range("B20").Value = "Geachte afzender,`nDank u voor uw e-mail. Helaas kunnen we geen bestellingen plaatsen via deze e-mail. Gelieve onze website te bezoeken en het bestelproces te doorlopen om uw M5-bouten te bestellen.`nMet vriendelijke groet,`n[E-mailassistent]"
$ws.Range("C20").Value = "Wil je 100 stuks M5-bouten bestellen?"
$ws.Range("D20").Value = "mailmind.test@zohomail.eu"
$ws.Range("E20").Value = "Bestelling / Levering"
$ws.Range("F20").Value = "2025-07-29 21:35:19"
$ws.Range("G20").Value = "Ja"
$ws.Range("H20").Value = "Nee"
$ws.Range("I20").Value = "Ja"
$ws.Range("J20").Value = "Nee"
